$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LMSData")
$ws.Activate()

# Update the Security Domains value for row 3 (J column) with new value
# (leading apostrophe preserves the existing "stored as text" quote-prefix
# cell style, matching the original formatting of this column)
$ws.Range("J3").Value = "'CORE TEST A;CORE TEST E"

# Update the active selection to I9 as recorded in the saved view state
$ws.Range("I9").Select()
